{"js": "// Update the default disk-space allowance mentioned in the Instant Access\n// form from 500 GiB to 250 GiB. Only the two sentences that describe the\n// *default* allowance change; the \"[Specify disk space requirements if\n// larger than 500 GiB]\" placeholder is left untouched.\n\n// 1) \"By default, projects are granted 500 GiB of disk space. ...\"\nconst grantedResults = context.document.body.search(\n  \"By default, projects are granted 500 GiB of disk space\",\n  { matchCase: true, matchWholeWord: false }\n);\ngrantedResults.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < grantedResults.items.length; i++) {\n  grantedResults.items[i].insertText(\n    \"By default, projects are granted 250 GiB of disk space\",\n    Word.InsertLocation.replace\n  );\n}\n\n// 2) \"Justification of the disk space requested if more than 500 GiB ...\"\nconst justifyResults = context.document.body.search(\n  \"disk space requested if more than 500 GiB\",\n  { matchCase: true, matchWholeWord: false }\n);\njustifyResults.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < justifyResults.items.length; i++) {\n  justifyResults.items[i].insertText(\n    \"disk space requested if more than 250 GiB\",\n    Word.InsertLocation.replace\n  );\n}\n\nawait context.sync();\n", "ps1": "# Update the default disk-space allowance mentioned in the Instant Access\n# form from 500 GiB to 250 GiB. Only the two sentences that describe the\n# *default* allowance change; the \"[Specify disk space requirements if\n# larger than 500 GiB]\" placeholder is left untouched.\n\n$d = $word.ActiveDocument\n\n# 1) \"By default, projects are granted 500 GiB of disk space. ...\"\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"granted 500 GiB of disk space\"\n$find.Replacement.Text = \"granted 250 GiB of disk space\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n# 2) \"Justification of the disk space requested if more than 500 GiB ...\"\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = \"disk space requested if more than 500 GiB\"\n$find2.Replacement.Text = \"disk space requested if more than 250 GiB\"\n$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2) | Out-Null\n"}
